$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.739.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.722.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5730"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.94%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "22.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06578"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07519"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.730.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.665"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5968"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.957.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "74.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("E17").Value = "  -9.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.700.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.271"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "203.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.553"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.991"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1219"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06186"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.378"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.387"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.699"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.699"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.671"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.029"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6425"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("E37").Value = "  -4.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.678"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01658"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.117.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.146"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8739"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.870.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000107"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.556"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.176"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05368"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4409"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.90%  "
